$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyIdentifiers")

# Make this the active sheet (matches activeTab moving from soa -> studyIdentifiers)
$ws.Activate()

# New "address" column (F) with header + sample pipe-delimited value
$ws.Range("F1").Value = "address"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = "line|city|district|state|postal_code|GBR"

# Match the widened column F seen in the target workbook
$ws.Columns.Item(6).ColumnWidth = 46.33

# Leave the selection sitting on the newly-added cell
$ws.Range("F2").Select()
